$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A24").Value = "nerdtree"
$ws.Range("B24").Value = "basic"
$ws.Range("C24").Value = "# Basic Ops`r`no: open files`r`nt: open in tab`r`ni: open in split`r`nu: upper level of directory`r`ne: explore selected direcotry`r`no/x: (de)select direcotry`r`n# How to open `r`n:NERDTree {direcotry}"

$ws.Range("C24").WrapText = $true

$ws.Rows.Item(24).RowHeight = 135

$ws.Range("C25").Select()
